$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '41.993.97'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.263.64'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.34%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '305.09'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '95.72'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.30%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.527'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.76%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('E9').Value = '  +0.45%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '35.01'
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0791'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.46%  '
$ws.Range('E12').Value = '  -0.41%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.76'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.16%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.612.23'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.42%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.39'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.50%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.263.85'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.60%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.791'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.75%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '41.856.46'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.39%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.40'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.18%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0902'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.85%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.96'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.41%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '67.99'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.32%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '237.35'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.80%  '
$ws.Range('E24').Value = '  -1.77%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.93'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.75%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '23.70'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.32%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '36.66'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +4.06%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.12'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.04%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '9.45'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.52%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '160.06'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.26%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.21'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.21%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('E34').Value = '  +4.74%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0735'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.99%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '17.06'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.64%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.82'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.26%  '
$ws.Range('E40').Value = '  -2.21%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.99'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.54%  '
$ws.Range('E42').Value = '  +1.98%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.970.89'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.13%  '
$ws.Range('E44').Value = '  -0.45%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '18.60'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -5.68%  '
$ws.Range('E46').Value = '  +0.68%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.86'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.75%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '52.92'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '72.49'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.27%  '
$ws.Range('E50').Value = '  -1.09%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '91.27'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.59%  '
